# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update row 2 (Home) target depth totals ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 203
$wsOff.Range("C2").Value = 126
$wsOff.Range("D2").Value = 43
$wsOff.Range("E2").Value = 16
$wsOff.Range("F2").Value = 4

# --- DEF sheet: update row 2 (Home) target depth totals ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 268
$wsDef.Range("C2").Value = 188
$wsDef.Range("D2").Value = 77
$wsDef.Range("E2").Value = 33
$wsDef.Range("G2").Value = 4
